# Update the "Year Pivot" table with refreshed detection-rate counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 3-11, columns C:H (B/"Years" column is unchanged).
$data = @{
    3  = @(64, 12, 15, 11, 7, 1491)
    4  = @(79, 18, 9, 17, 10, 1467)
    5  = @(80, 18, 11, 13, 11, 1467)
    6  = @(74, 23, 18, 10, 14, 1461)
    7  = @(82, 17, 11, 6, 10, 1474)
    8  = @(58, 14, 6, 14, 12, 1496)
    9  = @(36, 11, 1, 5, 7, 1540)
    10 = @(22, 4, 3, 6, 5, 1560)
    11 = @(12, 1, 3, 2, 4, 1578)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 3 + $i  # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Move the active selection from F7 to F6, matching the author's final state.
$ws.Range("F6").Select()
